$p = $ppt.ActivePresentation

# Slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS") contains a table
# (Google Shape;122;p17) that had a custom table style applied. Re-style
# it by applying a different (built-in) table style, identified by its
# style GUID, matching the Table Styles gallery action in PowerPoint's
# Table Design ribbon tab.
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)

if ($tableShape.HasTable) {
    $tableShape.Table.ApplyStyle("{2E41A223-C100-4A40-81A1-1614161B0F24}")
}
